$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh the "last updated" timestamp shown in the title cell ---------
$ws.Range("A1").Value = "Datos actualizados a 19 de Mayo de 2020 a las 07:35"

# --- Update the day's case counts for the countries whose totals moved ----
# These two countries received new figures, which pushes them past their
# neighbours once the table is re-ranked by "Casos totales" (column B).

# Hungria: new totals
$ws.Range("B69").Value = 3556
$ws.Range("C69").Value = 21
$ws.Range("D69").Value = 1412
$ws.Range("E69").Value = 1677
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 5
$ws.Range("H69").Value = 467

# Uzbekistan: new totals
$ws.Range("B75").Value = 2802
$ws.Range("C75").Value = 11
$ws.Range("D75").Value = 2314
$ws.Range("E75").Value = 475
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 13

# --- Re-rank the rows whose order changes as a result --------------------
# Country table is kept sorted descending by "Casos totales" (column B).
# Helper to (re)write a whole data row in one go.
function Set-CountryRow($row, $country, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Range("A$row").Value = $country
    $ws.Range("B$row").Value = $b
    $ws.Range("C$row").Value = $c
    $ws.Range("D$row").Value = $d
    $ws.Range("E$row").Value = $e
    $ws.Range("F$row").Value = $f
    $ws.Range("G$row").Value = $g
    $ws.Range("H$row").Value = $h
}

# Hungria now outranks Irak (row 69 <-> 70).
Set-CountryRow 69 "Hungria"    3556 21  1412 1677 0 5 467
Set-CountryRow 70 "Irak"       3554 0   2310 1117 0 0 127

# Uzbekistan now outranks Honduras & Guinea (rows 75-77 shift down one).
Set-CountryRow 75 "Uzbekistan" 2802 11  2314 475  0 0 13
Set-CountryRow 76 "Honduras"   2798 152 340  2312 0 4 146
Set-CountryRow 77 "Guinea"     2796 0   1263 1517 0 0 16

# Same day's update also re-ordered a few fully/partially tied groups
# further down the table (rows 195-197, 209-211, 214-216).
Set-CountryRow 195 "Santa Lucia"      18 0 18 0 0 0 0
Set-CountryRow 196 "Nueva Caledonia"  18 0 18 0 0 0 0
Set-CountryRow 197 "Belice"           18 0 16 0 0 0 2

Set-CountryRow 209 "Seychelles"   11 0 11 0 0 0 0
Set-CountryRow 210 "Groenlandia"  11 0 11 0 0 0 0
Set-CountryRow 211 "Montserrat"   11 0 10 0 0 0 1

Set-CountryRow 214 "Sahara Occidental"           6 0 6 0 0 0 0
Set-CountryRow 215 "San Bartolome"                6 0 6 0 0 0 0
Set-CountryRow 216 "Bonaire, San Eustaquio y Saba" 6 0 6 0 0 0 0
